# Apply weekly refresh of Fruta/Hortaliza data: the rows 2-32 of the sheet are
# re-populated by shuffling the data-bearing columns (D, H, I, J, K, L, M, O, P)
# among themselves according to the mapping below. Columns A, B, C, E, F, G, N,
# Q, R are constant across all rows and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: destination row -> source row (values are copied from the source row's
# "before" state into the destination row).
$rowMap = @{
    2  = 28
    3  = 23
    4  = 31
    5  = 30
    6  = 10
    7  = 11
    8  = 9
    9  = 20
    10 = 6
    11 = 2
    12 = 24
    13 = 19
    14 = 25
    15 = 32
    16 = 27
    17 = 18
    18 = 26
    19 = 8
    20 = 14
    21 = 29
    22 = 13
    23 = 21
    24 = 7
    25 = 5
    26 = 22
    27 = 16
    28 = 12
    29 = 17
    30 = 15
    31 = 3
    32 = 4
}

# Columns (by index) that get shuffled between rows.
# 4=D(Fecha) 8=H(Variedad) 9=I(Calidad) 10=J(Volumen) 11=K(Precio minimo)
# 12=L(Precio maximo) 13=M(Precio promedio ponderado) 15=O(Origen) 16=P(Precio $/Kg)
$cols = @(4, 8, 9, 10, 11, 12, 13, 15, 16)

# First, snapshot all the "before" values for every row/col we might touch,
# so that writes to earlier rows don't corrupt reads for later rows.
$snapshot = @{}
for ($r = 2; $r -le 32; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write the shuffled values into each destination row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
